$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) New "Rg / Rc / K / ec" mini-table next to the existing "Datos" block
#    (rows 2-5, columns E:F)
# ---------------------------------------------------------------------------
$ws.Range("E2").Value = "Rg [Ohm]"
$ws.Range("F2").Value = 120

$ws.Range("E3").Value = "Rc [Ohm]"
$ws.Range("F3").Value = 29880

$ws.Range("E4").Value = "K"
$ws.Range("F4").Value = 2.1

$ws.Range("E5").Value = "ec"
$ws.Range("F5").NumberFormat = "0.0000E+00"
$ws.Range("F5").Formula = "=F2/(F2+F3)*(1/F4)"

# ---------------------------------------------------------------------------
# 2) Un-merge & clear the old right-hand "Transversal abajo" /
#    "Longitudinal arriba y abajo" blocks (K9:Q16) - they get rebuilt
#    further down the sheet (rows 21-28).
# ---------------------------------------------------------------------------
$ws.Range("K9:M9").UnMerge()
$ws.Range("O9:Q9").UnMerge()
$ws.Range("K9:Q16").Clear()

# ---------------------------------------------------------------------------
# 3) Rework the first data table (rows 13-16): it gains a "Masa [kg]" /
#    "Vs [mV]" / "e [adim]" column set, while the old F/G (Tensión Vs /
#    Tensión MPa) columns move one slot to the right (E->F) losing the
#    "Vs" header and re-using it for the stress results.
# ---------------------------------------------------------------------------
$ws.Range("B13").Value = "Masa [kg]"
$ws.Range("C13").Value = "Carga [N]"
$ws.Range("D13").Value = "Vs [mV]"
$ws.Range("D13").Font.Bold = $true
$ws.Range("E13").Value = "e [adim]"
$ws.Range("E13").Font.Bold = $true
$ws.Range("F13").Value = "Tensión [MPa]"

# Clear the now-empty G13 header cell (old "Tensión [MPa]" header slot)
$ws.Range("G13").Clear()

# Row 14
$ws.Range("D14").Value = 280
$ws.Range("F14").NumberFormat = "0.00"
$ws.Range("F14").Formula = "=0.000001*6*`$C`$3*C14/(`$C`$4*`$C`$5*`$C`$5)"
$ws.Range("E14").NumberFormat = "0.0000E+00"
$ws.Range("E14").Formula = "=`$F`$5*D14/`$C`$11"
$ws.Range("G14").Clear()

# Row 15
$ws.Range("D15").Value = 560
$ws.Range("F15").NumberFormat = "0.00"
$ws.Range("F15").Formula = "=0.000001*6*`$C`$3*C15/(`$C`$4*`$C`$5*`$C`$5)"
$ws.Range("E15").NumberFormat = "0.0000E+00"
$ws.Range("E15").Formula = "=`$F`$5*D15/`$C`$11"
$ws.Range("G15").Clear()

# Row 16
$ws.Range("D16").Value = 840
$ws.Range("F16").NumberFormat = "0.00"
$ws.Range("F16").Formula = "=0.000001*6*`$C`$3*C16/(`$C`$4*`$C`$5*`$C`$5)"
$ws.Range("E16").NumberFormat = "0.0000E+00"
$ws.Range("E16").Formula = "=`$F`$5*D16/`$C`$11"
$ws.Range("G16").Clear()

# ---------------------------------------------------------------------------
# 4) Rebuild the "Transversal abajo" / "Longitudinal arriba y abajo" blocks
#    further down the sheet, at rows 21-28, columns B:D and F:H.
# ---------------------------------------------------------------------------
$ws.Range("B21:D21").Merge()
$ws.Range("B21:D21").Font.Bold = $true
$ws.Range("B21:D21").HorizontalAlignment = -4108
$ws.Range("B21").Value = "Transversal abajo"

$ws.Range("F21:H21").Merge()
$ws.Range("F21:H21").Font.Bold = $false
$ws.Range("F21:H21").HorizontalAlignment = -4108
$ws.Range("F21").Value = "Longitudinal arriba y abajo"

$ws.Range("B23").Value = "Vc [mV]"
$ws.Range("B23").Font.Bold = $true
$ws.Range("C23").Value = 1807

$ws.Range("F23").Value = "Vc [mV]"
$ws.Range("G23").Value = 1801

$ws.Range("B25").Value = "Peso [kg]"
$ws.Range("B25").Font.Bold = $true
$ws.Range("C25").Value = "Carga [N]"
$ws.Range("C25").Font.Bold = $true
$ws.Range("D25").Value = "Tensión [mV]"
$ws.Range("D25").Font.Bold = $true

$ws.Range("F25").Value = "Peso [kg]"
$ws.Range("F25").Font.Bold = $true
$ws.Range("G25").Value = "Carga [N]"
$ws.Range("G25").Font.Bold = $true
$ws.Range("H25").Value = "Tensión [mV]"
$ws.Range("H25").Font.Bold = $true

# Row 26
$ws.Range("B26").Value = 1
$ws.Range("D26").Value = 86
$ws.Range("F26").Value = 1
$ws.Range("G26").Formula = "=9.81*F26"
$ws.Range("H26").Value = 550

# Row 27
$ws.Range("B27").Value = 2
$ws.Range("D27").Value = 174
$ws.Range("F27").Value = 2
$ws.Range("G27").Formula = "=9.81*F27"
$ws.Range("H27").Value = 1110

# Row 28
$ws.Range("B28").Value = 3
$ws.Range("D28").Value = 260
$ws.Range("F28").Value = 3
$ws.Range("G28").Formula = "=9.81*F28"
$ws.Range("H28").Value = 1690

# ---------------------------------------------------------------------------
# 5) Final selection, matching the saved workbook view.
# ---------------------------------------------------------------------------
$ws.Range("F16").Select() | Out-Null
